# Apply "Updated LCA process names" edit to the Impact Assessment B - GWP workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the (only) worksheet: Sheet2 -> Sheet9
$ws.Name = "Sheet9"

# 2) Update the run date/time stamp in D1 (date, dd-MM-yyyy) and F1 (time, hh:mm)
#    45572 / 0.806519050925926  ->  45574 / 0.882288229166667
$ws.Range("D1").Value = 45574
$ws.Range("F1").Value = 0.882288229166667

# 3) Simplify the per-operation LCA process names in the header row (row 16),
#    columns E through U, to the new short naming scheme.
$ws.Range("E16").Value = "Turning 1"
$ws.Range("F16").Value = "Turning 2"
$ws.Range("G16").Value = "Turning 3"
$ws.Range("H16").Value = "Turning 4"
$ws.Range("I16").Value = "Turning 5"
$ws.Range("J16").Value = "Turning 6"
$ws.Range("K16").Value = "Turning 7"
$ws.Range("L16").Value = "Turning 8"
$ws.Range("M16").Value = "Turning 9"
$ws.Range("N16").Value = "Turning 10"
$ws.Range("O16").Value = "Turning 11"
$ws.Range("P16").Value = "Turning 13b"
$ws.Range("Q16").Value = "Drilling"
$ws.Range("R16").Value = "Milling"
$ws.Range("S16").Value = "Turning 14"
$ws.Range("T16").Value = "Surface Grinding"
$ws.Range("U16").Value = "Induction Hardening"
